$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1.47
$ws.Range("G2").Value = 1.48
$ws.Range("I2").Value = 8.2
$ws.Range("S2").Value = 2.64
$ws.Range("U2").Value = 2.12
$ws.Range("Z2").Value = 990
$ws.Range("F3").Value = 1.97
$ws.Range("G3").Value = 1.99
$ws.Range("H3").Value = 3.85
$ws.Range("I3").Value = 4
$ws.Range("K3").Value = 4.3
$ws.Range("M3").Value = 1.04
$ws.Range("P3").Value = 2.38
$ws.Range("S3").Value = 2.72
$ws.Range("U3").Value = 2.46
$ws.Range("Y3").Value = 19
$ws.Range("Z3").Value = 50
$ws.Range("AA3").Value = 1000
$ws.Range("AB3").Value = 12
$ws.Range("AD3").Value = 16.5
$ws.Range("AE3").Value = 1000
$ws.Range("AF3").Value = 14.5
$ws.Range("AH3").Value = 17
$ws.Range("AJ3").Value = 24
$ws.Range("AK3").Value = 19
$ws.Range("AL3").Value = 29
$ws.Range("AN3").Value = 10.5
$ws.Range("AO3").Value = 60
$ws.Range("H4").Value = 2.32
$ws.Range("K4").Value = 6.6
$ws.Range("P5").Value = 2.12
$ws.Range("Q5").Value = 1.7
$ws.Range("F6").Value = 2.48
$ws.Range("I6").Value = 3.4
$ws.Range("J6").Value = 3.35
$ws.Range("T6").Value = 1.87
$ws.Range("AJ6").Value = 36
$ws.Range("AK6").Value = 29
$ws.Range("G7").Value = 2.64
$ws.Range("H7").Value = 2.64
$ws.Range("I7").Value = 2.68
$ws.Range("N7").Value = 6.2
$ws.Range("P7").Value = 2.8
$ws.Range("R7").Value = 1.73
$ws.Range("T7").Value = 1.5
$ws.Range("AO7").Value = 13
$ws.Range("F8").Value = 9.6
$ws.Range("G8").Value = 10
$ws.Range("I8").Value = 1.39
$ws.Range("J8").Value = 5.7
$ws.Range("K8").Value = 5.8
$ws.Range("S8").Value = 2.6
$ws.Range("U8").Value = 1.97
$ws.Range("AB8").Value = 36
$ws.Range("AG8").Value = 36
$ws.Range("AI8").Value = 36
$ws.Range("AJ8").Value = 390
$ws.Range("AN8").Value = 200
$ws.Range("F9").Value = 1.19
$ws.Range("G9").Value = 1.2
$ws.Range("I9").Value = 22
$ws.Range("J9").Value = 8.2
$ws.Range("P9").Value = 2.86
$ws.Range("U9").Value = 1.71
$ws.Range("X9").Value = 34
$ws.Range("AE9").Value = 430
$ws.Range("AH9").Value = 110
$ws.Range("F10").Value = 1.29
$ws.Range("I10").Value = 11.5
$ws.Range("J10").Value = 7.2
$ws.Range("K10").Value = 7.4
$ws.Range("N10").Value = 8.4
$ws.Range("Q10").Value = 1.38
$ws.Range("R10").Value = 2
$ws.Range("U10").Value = 2.26
$ws.Range("AB10").Value = 15.5
$ws.Range("AC10").Value = 18
$ws.Range("AG10").Value = 11.5
$ws.Range("AH10").Value = 26
$ws.Range("F11").Value = 6.4
$ws.Range("G11").Value = 6.6
$ws.Range("J11").Value = 4.3
$ws.Range("K11").Value = 4.5
$ws.Range("N11").Value = 4.3
$ws.Range("Q11").Value = 1.83
$ws.Range("R11").Value = 1.44
$ws.Range("T11").Value = 1.91
$ws.Range("U11").Value = 2.02
$ws.Range("Y11").Value = 9
$ws.Range("AA11").Value = 15
$ws.Range("AC11").Value = 9.8
$ws.Range("AD11").Value = 9.8
$ws.Range("AO11").Value = 8.2
$ws.Range("H12").Value = 2.42
$ws.Range("I12").Value = 2.44
$ws.Range("S12").Value = 2.9
$ws.Range("X12").Value = 19
$ws.Range("AA12").Value = 34
$ws.Range("AE12").Value = 24
$ws.Range("AN12").Value = 24
$ws.Range("F13").Value = 2.48
$ws.Range("G13").Value = 2.84
$ws.Range("H13").Value = 2.84
$ws.Range("I13").Value = 3.65
$ws.Range("J13").Value = 3.1
$ws.Range("K13").Value = 3.7
$ws.Range("P13").Value = 1.78
$ws.Range("Q13").Value = 2.04
